# Fix issue for name and unit values not appearing in iModel
$wb = $excel.ActiveWorkbook

$wsDevice = $wb.Worksheets.Item(1)             # "Device"
$wsTemperature = $wb.Worksheets.Item(2)        # "TemperatureDatapoint"
$wsPressure = $wb.Worksheets.Item(3)           # "PressureDatapoint"

# --- Device sheet: row 6 was miscategorized as a Pressure device (mm hg) with
#     id/name "7" / "Device 7"; it should be device 8, a Temperature device (Celsius).
$wsDevice.Range("A6").Value = 8
$wsDevice.Range("B6").Value = "Device 8"
$wsDevice.Range("C6").Value = "Temperature "
$wsDevice.Range("D6").Value = "Celsius"

# --- PressureDatapoint sheet: correct the sample values in column C.
$wsPressure.Range("C2").Value = 4
$wsPressure.Range("C3").Value = 4

# --- Re-point the active sheet / selections back to "Device" (tab 1),
#     leaving the other two sheets with updated selections but not activated.
$wsTemperature.Activate()
$wsTemperature.Range("B2").Select()

$wsPressure.Activate()
$wsPressure.Range("C3").Select()

$wsDevice.Activate()
$wsDevice.Range("D6").Select()
